$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-converted to numbers by Excel (losing formatting like trailing zeros).
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '47.840.51'
$ws.Range("E2").Value = '  +5.92%  '
$ws.Range("D3").Value = '2.513.96'
$ws.Range("E3").Value = '  +3.47%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '324.21'
$ws.Range("E5").Value = '  +2.26%  '
$ws.Range("D6").Value = '106.30'
$ws.Range("E6").Value = '  +3.43%  '
$ws.Range("E7").Value = '  +1.76%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").Value = '0.543'
$ws.Range("E9").Value = '  +3.04%  '
$ws.Range("D10").Value = '38.00'
$ws.Range("E10").Value = '  +7.02%  '
$ws.Range("D11").Value = '0.0817'
$ws.Range("E11").Value = '  +1.71%  '
$ws.Range("E12").Value = '  +0.86%  '
$ws.Range("D13").Value = '18.43'
$ws.Range("E13").Value = '  +1.53%  '
$ws.Range("E14").Value = '  +2.17%  '
$ws.Range("D15").Value = '2.908.80'
$ws.Range("E15").Value = '  +3.56%  '
$ws.Range("D16").Value = '2.527.96'
$ws.Range("E16").Value = '  +3.82%  '
$ws.Range("D17").Value = '0.849'
$ws.Range("E17").Value = '  +0.87%  '
$ws.Range("D18").Value = '47.728.77'
$ws.Range("E18").Value = '  +5.84%  '
$ws.Range("D19").Value = '12.75'
$ws.Range("E19").Value = '  +3.68%  '
$ws.Range("D20").Value = '6.58'
$ws.Range("E20").Value = '  +3.36%  '
$ws.Range("D21").Value = '0.0₃0939'
$ws.Range("E21").Value = '  +1.90%  '
$ws.Range("D22").Value = '70.85'
$ws.Range("E22").Value = '  +2.89%  '
$ws.Range("D23").Value = '251.61'
$ws.Range("E23").Value = '  +3.03%  '
$ws.Range("E24").Value = '  +6.60%  '
$ws.Range("D25").Value = '2.57'
$ws.Range("E25").Value = '  +3.04%  '
$ws.Range("D26").Value = '26.38'
$ws.Range("E26").Value = '  +3.36%  '
$ws.Range("E27").Value = '  -0.10%  '
$ws.Range("E28").Value = '  +4.96%  '
$ws.Range("D29").Value = '2.21'
$ws.Range("E29").Value = '  +6.83%  '
$ws.Range("E30").Value = '  +6.92%  '
$ws.Range("E31").Value = '  +9.59%  '
$ws.Range("E32").Value = '  +0.57%  '
$ws.Range("D33").Value = '20.11'
$ws.Range("E33").Value = '  -1.27%  '
$ws.Range("E35").Value = '  +2.72%  '
$ws.Range("E36").Value = '  +0.19%  '
$ws.Range("D37").Value = '1.95'
$ws.Range("E37").Value = '  +4.11%  '
$ws.Range("D38").Value = '4.64'
$ws.Range("E38").Value = '  +4.81%  '
$ws.Range("E39").Value = '  +4.69%  '
$ws.Range("B40").Value = 'WEMIXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D40").Value = '2.26'
$ws.Range("E40").Value = '  +2.24%  '
$ws.Range("B41").Value = 'Stellar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D41").Value = '0.112'
$ws.Range("E41").Value = '  +2.30%  '
$ws.Range("D42").Value = '121.79'
$ws.Range("D43").Value = '21.15'
$ws.Range("E43").Value = '  +2.70%  '
$ws.Range("D44").Value = '0.0298'
$ws.Range("E44").Value = '  +3.46%  '
$ws.Range("D45").Value = '1.972.21'
$ws.Range("E45").Value = '  +1.91%  '
$ws.Range("D46").Value = '3.02'
$ws.Range("E46").Value = '  +3.33%  '
$ws.Range("E47").Value = '  -0.57%  '
$ws.Range("D48").Value = '1.82'
$ws.Range("E48").Value = '  +0.78%  '
$ws.Range("D49").Value = '9.21'
$ws.Range("E49").Value = '  -0.29%  '
$ws.Range("D50").Value = '5.38'
$ws.Range("E50").Value = '  +14.04%  '
$ws.Range("D51").Value = '79.26'
$ws.Range("E51").Value = '  +3.65%  '
